$d = $word.ActiveDocument

# Locate the paragraph that immediately follows the "LOM3038: ..." requisito
# paragraph (a blank paragraph), and the two trailing footer paragraphs that
# need to be removed: "Ver no Jupiter Salvar em pdf Salvar em docx" and the
# copyright/footer line. These three consecutive paragraphs are deleted,
# leaving the blank paragraph that sits just before the final page-break
# paragraph untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        # The blank paragraph right before this one is the start of the
        # block that must be removed.
        $startPara = $d.Paragraphs.Item($i - 1)
    }

    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
